# Updated Risk Register with owners of risks for risks related to a single WP
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risks")

# Clear the "New Probability" (G) and "Impact" (H) values for all risk rows (3-27).
# This resets the "Assessment" (I) formula results to their un-evaluated boolean state.
for ($r = 3; $r -le 27; $r++) {
    $ws.Range("G$r`:H$r").ClearContents()
}

# Assign an Owner (column K) to every risk that is related to a single Work Package.
# Risks tied to "All WPs" or to multiple WPs (e.g. "3,5") are left without an owner.
$owners = @{
    6  = "Andy Gotz"
    7  = "Andy Gotz"
    8  = "Tobias Richter"
    9  = "Tobias Richter"
    10 = "Hans Fangohr"
    11 = "Hans Fangohr"
    12 = "Carsten Fortmann-Grote"
    15 = "Carsten Fortmann-Grote"
    16 = "Jean-François Perrin"
    17 = "Jean-François Perrin"
    18 = "Roberto Pugliese"
    19 = "Roberto Pugliese"
    20 = "Roberto Pugliese"
    21 = "Roberto Pugliese"
    22 = "Thomas Rod"
    25 = "Nicoletta Carboni"
    26 = "Nicoletta Carboni"
    27 = "Nicoletta Carboni"
}

foreach ($r in $owners.Keys) {
    $ws.Range("K$r").Value = $owners[$r]
}

# Update the active selection on the sheet.
[void]$ws.Range("C3").Select()
